$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D..Y to E..Z).
$ws.Range("D1").EntireColumn.Insert()

# Excel normally carries the left-neighbour column's width onto a freshly
# inserted column; reproduce that explicitly (new col D == old col C width).
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Populate header + row-2 value for the new "brand_id" column.
$ws.Range("D1").Value = "brand_id"
$ws.Range("D2").Value = 1

# Match the authored selection state.
$ws.Range("D1:D2").Select() | Out-Null
